$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "62.208.91"
$ws.Range("E2").Value = "  -2.67%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.512.81"
$ws.Range("E3").Value = "  -4.08%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "557.01"
$ws.Range("E5").Value = "  -3.42%  "

# Row 6 - Solana
Set-TextValue "D6" "148.50"
$ws.Range("E6").Value = "  -5.20%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -3.03%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.515.51"
$ws.Range("E9").Value = "  -3.89%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -8.11%  "

# Row 11 - now Toncoin (was TRON)
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D11" "5.45"
$ws.Range("E11").Value = "  -6.49%  "

# Row 12 - now TRON (was Toncoin)
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D12" "0.155"
$ws.Range("E12").Value = "  -1.25%  "

# Row 13 - Cardano
Set-TextValue "D13" "0.365"
$ws.Range("E13").Value = "  -4.59%  "

# Row 14 - Avalanche
Set-TextValue "D14" "26.65"
$ws.Range("E14").Value = "  -5.88%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.966.71"
$ws.Range("E15").Value = "  -4.09%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -7.92%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "62.056.25"
$ws.Range("E17").Value = "  -2.60%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.492.56"
$ws.Range("E18").Value = "  -4.80%  "

# Row 19 - Chainlink
Set-TextValue "D19" "11.43"
$ws.Range("E19").Value = "  -5.73%  "

# Row 20 - Uniswap
Set-TextValue "D20" "7.17"
$ws.Range("E20").Value = "  -6.85%  "

# Row 21 - Polkadot
Set-TextValue "D21" "4.29"
$ws.Range("E21").Value = "  -6.10%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "326.17"
$ws.Range("E22").Value = "  -5.05%  "

# Row 23 - Dai
Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  -0.03%  "

# Row 24 - Litecoin
Set-TextValue "D24" "64.80"
$ws.Range("E24").Value = "  -4.13%  "

# Row 25 - SuiNetwork
Set-TextValue "D25" "1.77"
$ws.Range("E25").Value = "  -0.06%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -5.04%  "

# Row 27 - WrappedeETH
Set-TextValue "D27" "2.652.72"
$ws.Range("E27").Value = "  -3.30%  "

# Row 28 - now Fetch.AI (was InternetComputer(DFINITY))
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D28" "1.53"
$ws.Range("E28").Value = "  -2.88%  "

# Row 29 - now InternetComputer(DFINITY) (was Fetch.AI)
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D29" "8.67"
$ws.Range("E29").Value = "  -6.07%  "

# Row 30 - now Bittensor (was Binance-PegBSC-USD)
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D30" "544.24"
$ws.Range("E30").Value = "  -8.66%  "

# Row 31 - now Binance-PegBSC-USD (was Bittensor)
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D31" "0.999"
$ws.Range("E31").Value = "  -0.11%  "

# Row 32 - Aptos
Set-TextValue "D32" "7.83"
$ws.Range("E32").Value = "  -1.40%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  -4.11%  "

# Row 34 - PancakeSwap
Set-TextValue "D34" "1.93"
$ws.Range("E34").Value = "  -6.57%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -8.07%  "

# Row 36 - RenderToken
Set-TextValue "D36" "6.05"
$ws.Range("E36").Value = "  -8.72%  "

# Row 37 - NEARProtocol
Set-TextValue "D37" "4.97"
$ws.Range("E37").Value = "  -7.77%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  +0.12%  "

# Row 39 - PolygonEcosystemToken
Set-TextValue "D39" "0.385"
$ws.Range("E39").Value = "  -5.79%  "

# Row 40 - EthereumClassic
Set-TextValue "D40" "18.83"
$ws.Range("E40").Value = "  -4.85%  "

# Row 41 - Monero
Set-TextValue "D41" "151.51"

# Row 42 - Stacks
$ws.Range("E42").Value = "  -7.04%  "

# Row 43 - USDe
Set-TextValue "D43" "0.999"
$ws.Range("E43").Value = "  -0.04%  "

# Row 44 - dogwifhat
Set-TextValue "D44" "2.32"
$ws.Range("E44").Value = "  -4.42%  "

# Row 45 - Aave
Set-TextValue "D45" "151.27"
$ws.Range("E45").Value = "  -3.10%  "

# Row 46 - Filecoin
Set-TextValue "D46" "3.70"
$ws.Range("E46").Value = "  -5.73%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "22.06"
$ws.Range("E47").Value = "  -5.85%  "

# Row 48 - Hedera
Set-TextValue "D48" "0.0554"
$ws.Range("E48").Value = "  -6.65%  "

# Row 49 - Mantle
Set-TextValue "D49" "0.599"
$ws.Range("E49").Value = "  -4.92%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  -6.16%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -4.60%  "
